# JOB_FR_CONFIG.xlsx edit: merge-rename 3 source sheets and add SALAIRE / SECTEUR
# config rows to the JOB_FR sheet (matches TYPE_EMPLOI's row pattern).

$wb = $excel.ActiveWorkbook

# 1. Rename the three merged-source sheets.
$wb.Worksheets.Item("CGPME-FACEBOOK").Name = "FACEBOOKMERGED"
$wb.Worksheets.Item("MONSTER").Name = "MONSTERMERGED"
$wb.Worksheets.Item("POLEEMPLOI").Name = "POLEEMPLOIMERGED"

# 2. Add the two new config rows (SALAIRE, SECTEUR) to JOB_FR, mirroring the
#    existing TYPE_EMPLOI row (row 30): same rule name/parameters/logs.
$jobFr = $wb.Worksheets.Item("JOB_FR")

$jobFr.Range("A31").Value = "expect_column_values_to_not_be_null"
$jobFr.Range("B31").Value = "SALAIRE"
$jobFr.Range("C31").Value = "{'mostly': 0.9}"
$jobFr.Range("D31").Value = "{""keyword"": ""FILE_CHECK""}"

$jobFr.Range("A32").Value = "expect_column_values_to_not_be_null"
$jobFr.Range("B32").Value = "SECTEUR"
$jobFr.Range("C32").Value = "{'mostly': 0.9}"
$jobFr.Range("D32").Value = "{""keyword"": ""FILE_CHECK""}"

# 3. Restore the saved cursor/selection positions.
$jobFr.Range("C36").Select()

$facebookMerged = $wb.Worksheets.Item("FACEBOOKMERGED")
$facebookMerged.Range("J29").Select()

$jobFr.Select()
